# Append a new row of Argent (solar) price data to the "Prices" sheet.
# Row 80 mirrors the layout of the existing rows (all values stored as
# literal text, same as the rest of the data table), so each value is
# written with a leading apostrophe to force Excel to keep it as text
# instead of auto-converting to a date/number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 80

$ws.Cells.Item($row, 1).Value = "'2025-05-20"
$ws.Cells.Item($row, 2).Value = "'37"
$ws.Cells.Item($row, 3).Value = "'37"
$ws.Cells.Item($row, 4).Value = "'0.94"
$ws.Cells.Item($row, 5).Value = "'0.258"
$ws.Cells.Item($row, 6).Value = "'0.09"
$ws.Cells.Item($row, 7).Value = "'5,299"
$ws.Cells.Item($row, 8).Value = "'7,933"
$ws.Cells.Item($row, 9).Value = "'7,983"
$ws.Cells.Item($row, 10).Value = "'7.2266"
